$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '65.438.21'
$ws.Range("E2").Value = '  -2.40%  '
$ws.Range("D3").Value = '3.481.83'
$ws.Range("E3").Value = '  +0.44%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.00'
$ws.Range("E4").Value = '  +0.16%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '553.10'
$ws.Range("E5").Value = '  +1.18%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '178.63'
$ws.Range("E6").Value = '  -3.99%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.638'
$ws.Range("E7").Value = '  +4.38%  '
$ws.Range("E8").Value = '  +0.06%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.632'
$ws.Range("E9").Value = '  -0.74%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.151'
$ws.Range("E10").Value = '  +2.74%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '53.65'
$ws.Range("E11").Value = '  -4.17%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.0000270'
$ws.Range("E12").Value = '  -1.41%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '9.22'
$ws.Range("E13").Value = '  -1.99%  '
$ws.Range("D14").Value = '4.039.15'
$ws.Range("E14").Value = '  +0.59%  '
$ws.Range("E15").Value = '  +2.22%  '
$ws.Range("D16").Value = '3.480.17'
$ws.Range("E16").Value = '  +0.44%  '
$ws.Range("E17").Value = '  +0.43%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '12.03'
$ws.Range("E18").Value = '  +2.16%  '
$ws.Range("D19").Value = '65.490.99'
$ws.Range("E19").Value = '  -2.90%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.988'
$ws.Range("E20").Value = '  -1.79%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '416.05'
$ws.Range("E21").Value = '  +2.65%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '4.04'
$ws.Range("E22").Value = '  +3.46%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '86.10'
$ws.Range("E23").Value = '  +1.74%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '4.11'
$ws.Range("E24").Value = '  -1.94%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '12.96'
$ws.Range("E25").Value = '  +10.13%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '10.76'
$ws.Range("E26").Value = '  -10.11%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '2.84'
$ws.Range("E27").Value = '  -2.19%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '6.04'
$ws.Range("E28").Value = '  -3.60%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '9.09'
$ws.Range("E29").Value = '  +5.43%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '30.23'
$ws.Range("E30").Value = '  +0.34%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '6.52'
$ws.Range("E31").Value = '  -4.90%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '612.14'
$ws.Range("E32").Value = '  -9.85%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '11.74'
$ws.Range("E33").Value = '  +0.91%  '
$ws.Range("E34").Value = '  -0.07%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '59.20'
$ws.Range("E35").Value = '  +0.33%  '
$ws.Range("B36").Value = 'Kaspa'
$ws.Range("C36").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.146'
$ws.Range("E36").Value = '  +10.09%  '
$ws.Range("B37").Value = 'Dai'
$ws.Range("C37").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.00'
$ws.Range("E37").Value = '  +0.19%  '
$ws.Range("B38").Value = 'InjectiveProtocol'
$ws.Range("C38").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '37.37'
$ws.Range("E38").Value = '  -2.71%  '
$ws.Range("B39").Value = 'Maker'
$ws.Range("C39").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D39").Value = '3.377.83'
$ws.Range("E39").Value = '  +11.08%  '
$ws.Range("B40").Value = 'PEPE'
$ws.Range("C40").Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range("D40").Value = '0.0₃0782'
$ws.Range("E40").Value = '  -5.39%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.379'
$ws.Range("E41").Value = '  -5.84%  '
$ws.Range("E42").Value = '  -0.07%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '3.24'
$ws.Range("E43").Value = '  -4.64%  '
$ws.Range("E44").Value = '  -5.57%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '2.53'
$ws.Range("E45").Value = '  -9.51%  '
$ws.Range("B46").Value = 'VeChain'
$ws.Range("C46").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.0414'
$ws.Range("E46").Value = '  -1.30%  '
$ws.Range("B47").Value = 'ApeXProtocol'
$ws.Range("C47").Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '3.25'
$ws.Range("E47").Value = '  -0.53%  '
$ws.Range("E48").Value = '  -1.17%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.132'
$ws.Range("E49").Value = '  +2.19%  '
$ws.Range("E50").Value = '  -4.00%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '137.80'
$ws.Range("E51").Value = '  -1.53%  '
